$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shaconemo ping-file revision bump: r255 -> r256 across every
# "Identified in the shaconemo (r255) ... ping file: ..." comment cell.
$rng = $ws.UsedRange
$null = $rng.Replace("(r255)", "(r256)", -4163, 1, $false, $false, $false)

# A handful of seaIce comments also got their underlying ping-file
# variable names updated to include an "_sc_" infix.
$null = $rng.Replace("NH_icearea", "NH_sc_icearea", -4163, 1, $false, $false, $false)
$null = $rng.Replace("SH_icearea", "SH_sc_icearea", -4163, 1, $false, $false, $false)
$null = $rng.Replace("NH_iceextt", "NH_sc_iceextt", -4163, 1, $false, $false, $false)
$null = $rng.Replace("SH_iceextt", "SH_sc_iceextt", -4163, 1, $false, $false, $false)
$null = $rng.Replace("NH_icevolu", "NH_sc_icevolu", -4163, 1, $false, $false, $false)
$null = $rng.Replace("SH_icevolu", "SH_sc_icevolu", -4163, 1, $false, $false, $false)

# Update the sheet's saved view/selection: was scrolled to show column G
# with G3 selected; now scrolled back to A1 with C3:C269 selected.
$null = $ws.Range("C3:C269").Select()

# Column A's stored width nudges from 10.66 to 10.65 characters.
$ws.Columns.Item(1).ColumnWidth = 9.75
